$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.273.75'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.143.76'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.39'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '149.98'
$ws.Range('E6').Value = '  -4.72%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.145.11'
$ws.Range('E8').Value = '  -2.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('E10').Value = '  -4.61%  '
$ws.Range('E11').Value = '  -1.40%  '
$ws.Range('E12').Value = '  -4.92%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000260'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '37.02'
$ws.Range('E14').Value = '  -4.37%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.620.79'
$ws.Range('E15').Value = '  -3.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.418.01'
$ws.Range('E16').Value = '  -3.11%  '
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.141.89'
$ws.Range('E18').Value = '  -2.33%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.05'
$ws.Range('E19').Value = '  -4.44%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '484.79'
$ws.Range('E20').Value = '  -5.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.73'
$ws.Range('E21').Value = '  -3.13%  '
$ws.Range('E22').Value = '  -2.45%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.81'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.89'
$ws.Range('E24').Value = '  -5.25%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.61'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -2.56%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.68'
$ws.Range('E28').Value = '  -4.72%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').Value = '  -4.27%  '
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.10'
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.72'
$ws.Range('E32').Value = '  -7.36%  '
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '26.88'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.11'
$ws.Range('E35').Value = '  -5.64%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.12'
$ws.Range('E36').Value = '  -6.23%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '54.65'
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.27'
$ws.Range('E38').Value = '  +6.33%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0₃0755'
$ws.Range('E39').Value = '  -3.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '451.53'
$ws.Range('E40').Value = '  -10.38%  '
$ws.Range('E41').Value = '  -3.76%  '
$ws.Range('E42').Value = '  -4.63%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.53'
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.893.38'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.41'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.275'
$ws.Range('E46').Value = '  -8.09%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.95'
$ws.Range('E47').Value = '  -4.73%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  -3.58%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '120.22'
$ws.Range('E51').Value = '  -1.53%  '
